$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (prices & percentage changes
# refreshed by the symbol-list update). Values must remain plain text
# (inlineStr in the source), matching the existing column formatting,
# so we force Text number format before assigning and clear it right
# after so no stray styling is left behind on the cell.
$updates = @{
    'D2' = '308.22'
    'E2' = '-1.33%'
    'D3' = '36.15'
    'E3' = '-3.88%'
    'E4' = '-0.82%'
    'D5' = '0.07700'
    'E5' = '-2.75%'
    'D6' = '4.392'
    'E6' = '-0.62%'
    'D7' = '8.297'
    'E7' = '0.36%'
    'D8' = '1.842'
    'E8' = '-3.41%'
    'E9' = '1.05%'
    'D10' = '0.9204'
    'E10' = '-0.03%'
    'D11' = '0.1100'
    'E11' = '-9.92%'
    'D12' = '0.1850'
    'E12' = '-4.05%'
    'D13' = '0.08740'
    'E13' = '-4.34%'
    'D14' = '0.03338'
    'E14' = '0.94%'
    'D15' = '0.09533'
    'E15' = '-0.86%'
    'D16' = '0.001383'
    'E16' = '0.21%'
    'D17' = '0.006093'
    'E17' = '4.56%'
    'E18' = '-4.59%'
    'D19' = '0.3436'
    'E19' = '-0.49%'
    'D20' = '6.328'
    'E20' = '20.13%'
    'D21' = '0.1291'
    'E21' = '1.43%'
    'E23' = '-1.06%'
    'D24' = '0.001200'
    'E24' = '-3.97%'
    'D25' = '0.004258'
    'E25' = '-1.17%'
    'D26' = '0.0001328'
    'E26' = '8.86%'
    'D27' = '0.0002902'
    'D39' = '0.02087'
    'E39' = '-2.87%'
    'D40' = '0.04936'
    'E40' = '-3.56%'
    'D41' = '0.007504'
    'E41' = '-0.77%'
    'D42' = '0.1352'
    'E42' = '-0.64%'
    'D43' = '0.008509'
    'E43' = '-5.25%'
    'D44' = '0.002067'
    'E44' = '2.83%'
    'D45' = '0.008383'
    'E45' = '-2.76%'
    'D46' = '0.00006318'
    'E46' = '-5.91%'
    'E47' = '0.07%'
    'D48' = '0.002847'
    'E48' = '-14.37%'
    'D49' = '0.001691'
    'E49' = '40.80%'
    'D50' = '0.00002101'
    'E50' = '0.07%'
    'D51' = '0.0002001'
    'E51' = '0.07%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}

Write-Output ("Updated " + $updates.Count + " cells")
